$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted with "." as thousands separators; force
# text storage so Excel does not reinterpret them as numbers/dates on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.623.76"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "1.844.38"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "316.15"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "0.4234"
$ws.Range("E7").Value = "  -2.75%  "

$ws.Range("D8").Value = "0.3637"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").Value = "45.23"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "0.07249"
$ws.Range("E10").Value = "  -3.22%  "

$ws.Range("D11").Value = "0.8907"
$ws.Range("E11").Value = "  -4.69%  "

$ws.Range("D12").Value = "20.61"
$ws.Range("E12").Value = "  -3.38%  "

$ws.Range("D13").Value = "1.882.17"
$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("D14").Value = "6.577"
$ws.Range("E14").Value = "  -1.76%  "

$ws.Range("D15").Value = "5.347"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "0.06867"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "78.71"
$ws.Range("E18").Value = "  -3.27%  "

$ws.Range("D19").Value = "0.000008845"
$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("D20").Value = "1.0000"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "15.45"
$ws.Range("E21").Value = "  -2.23%  "

$ws.Range("D22").Value = "27.609.50"
$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").Value = "4.981"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("D25").Value = "2.056.45"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "2.010"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").Value = "154.50"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  +1.03%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "119.44"
$ws.Range("E29").Value = "  +5.79%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.233"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").Value = "1.840"
$ws.Range("E31").Value = "  +6.34%  "

$ws.Range("D32").Value = "0.08906"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").Value = "0.7782"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").Value = "4.560"
$ws.Range("E34").Value = "  -5.19%  "

$ws.Range("D35").Value = "2.955"
$ws.Range("E35").Value = "  -2.53%  "

$ws.Range("D36").Value = "1.100"
$ws.Range("E36").Value = "  -6.10%  "

$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05393"
$ws.Range("E38").Value = "  -0.39%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.096"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").Value = "2.810"
$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("D42").Value = "6.872"
$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("D43").Value = "0.5061"
$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("D44").Value = "0.1646"
$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("D45").Value = "8.261"
$ws.Range("E45").Value = "  -5.23%  "

$ws.Range("D46").Value = "0.06606"
$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("D47").Value = "10.36"
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").Value = "0.4702"
$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").Value = "104.72"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").Value = "1.628"
$ws.Range("E51").Value = "  -2.35%  "
